# Add speaker notes to each slide of the presentation.
$p = $ppt.ActivePresentation

$notesText = @(
    "Speaker notes for slide 1: Introduction to contract update.",
    "Speaker notes for slide 2: Summary of key changes.",
    "Speaker notes for slide 3: Required steps for completion.",
    "Speaker notes for slide 4: Contact information and support."
)

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $notesPage = $slide.NotesPage
    $notesShape = $notesPage.Shapes.AddPlaceholder(2)
    $notesShape.TextFrame.TextRange.Text = $notesText[$i - 1]
}
